# Zeiterfassung_Jeckle_Lukas.xlsx - neue Zeitposten gebucht
#
# Adds 8 new time-tracking rows (31-38) below the existing log and moves
# the "Gesamt:" (total) row from row 32 down to row 40, widening its SUM
# formula to cover the new entries.
#
# Notes on ordering:
#  - The "Gesamt:" row is relocated *first*, before row 32 is reused for a
#    new entry, so that shared string is never orphaned mid-script.
#  - The brand-new description/label texts are written in the exact
#    sequence they need to be interned in (so the rebuilt shared-strings
#    table lines up with the target: new entries are appended in the order
#    they were first typed, not in final row order).
#  - Date/duration cells copy number formatting from an existing date/
#    duration cell of the same kind (xlPasteFormats) so Excel reuses the
#    existing style slots instead of minting new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Relocate the "Gesamt:" total row from 32 -> 40 first -----------------
$ws.Range("A32").Copy()
$ws.Range("A40").PasteSpecial($xlPasteFormats)
$ws.Range("B32").Copy()
$ws.Range("B40").PasteSpecial($xlPasteFormats)
$ws.Range("A40").Value = "Gesamt:"
$ws.Range("B40").Formula = "=SUM(B7:B38)"
$ws.Range("A32:B32").ClearContents()

# --- Copy date (A) / duration (B) column formatting onto the new rows ----
# (style "s=4" for dates, "s=2" for durations - reused from existing rows)
$ws.Range("A7").Copy(); $ws.Range("A31").PasteSpecial($xlPasteFormats)
$ws.Range("B7").Copy(); $ws.Range("B31").PasteSpecial($xlPasteFormats)

$ws.Range("A8").Copy(); $ws.Range("A32").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy(); $ws.Range("B32").PasteSpecial($xlPasteFormats)

$ws.Range("A7").Copy(); $ws.Range("A33").PasteSpecial($xlPasteFormats)
$ws.Range("B7").Copy(); $ws.Range("B33").PasteSpecial($xlPasteFormats)

$ws.Range("A8").Copy(); $ws.Range("A34").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy(); $ws.Range("B34").PasteSpecial($xlPasteFormats)

$ws.Range("A7").Copy(); $ws.Range("A35").PasteSpecial($xlPasteFormats)
$ws.Range("B7").Copy(); $ws.Range("B35").PasteSpecial($xlPasteFormats)

$ws.Range("A8").Copy(); $ws.Range("A36").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy(); $ws.Range("B36").PasteSpecial($xlPasteFormats)

$ws.Range("A8").Copy(); $ws.Range("A37").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy(); $ws.Range("B37").PasteSpecial($xlPasteFormats)

$ws.Range("A8").Copy(); $ws.Range("A38").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy(); $ws.Range("B38").PasteSpecial($xlPasteFormats)

# --- Write the brand-new "Buchungsposten"/"Beschreibung" text in the ------
# exact order they were first entered, so the rebuilt shared-strings table
# matches (existing entries like "Planung"/"Coding"/"Online-Meeting" are
# simply reused and don't affect ordering).
$ws.Range("C31").Value = "Präsenz-Meeting"
$ws.Range("D32").Value = "UML Diagramm grob umgesetzt."
$ws.Range("D36").Value = "Bugfixes implementiert."
$ws.Range("D38").Value = "Sprint Review #2"
$ws.Range("D31").Value = "Weekly Meeting #4 - Erste Version des Programm UML Diagramms erstellt."
$ws.Range("D33").Value = "Logo Quellen im Projekt hinterlegt & Möglichkeiten zur Verbesserung der Team-Zusammenarbeit im Projekt mittels GitHub Branches recherchiert."
$ws.Range("D35").Value = "Weekly Meeting #5 Summup erstellt."
$ws.Range("D34").Value = "Weekly Meeting #5 abgehalten."

# --- Remaining cells (dates, durations, re-used "Buchungsposten" labels) --
$ws.Range("A31").Value = 45608
$ws.Range("B31").Value = 2

$ws.Range("A32").Value = 45608
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "Coding"

$ws.Range("A33").Value = 45613
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "Planung"

$ws.Range("A34").Value = 45614
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = "Online-Meeting"

$ws.Range("A35").Value = 45614
$ws.Range("B35").Value = 0.5
$ws.Range("C35").Value = "Planung"

$ws.Range("A36").Value = 45614
$ws.Range("B36").Value = 1.5
$ws.Range("C36").Value = "Coding"

$ws.Range("A37").Value = 45616
$ws.Range("B37").Value = 1.5
$ws.Range("C37").Value = "Planung"
$ws.Range("D37").Value = "Aufgaben-, Projekt-, und Zeitmanagementplanung."

$ws.Range("A38").Value = 45617
$ws.Range("B38").Value = 0.5
$ws.Range("C38").Value = "Online-Meeting"

# --- Sheet-view tweaks: zoom 220% -> 160%, move the active selection -----
$excel.ActiveWindow.Zoom = 160
$ws.Range("B38").Select()
